$d = $word.ActiveDocument

# Make sure formatting/text edits are applied cleanly, without Word's
# "track changes" markup polluting the OOXML.
$d.TrackRevisions = $false

# ------------------------------------------------------------------
# 1. Merge the two runs around the stray mid-sentence "_GoBack" bookmark
#    in the Julia/Colab installation paragraph back into a single run.
# ------------------------------------------------------------------
$oldText = "cell needs to be run in order to install Julia, whenever you restart your browser. Please follow the instructions in the "
$rngMerge = $d.Content
$rngMerge.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

# ------------------------------------------------------------------
# 2. Shrink the "Installation Guide" heading from 14pt to 13pt
#    (sz/szCs 28 -> 26) and re-anchor the "_GoBack" bookmark there.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$rngTitle = $p2.Range
$rngTitle.Font.Size = 13
$rngTitle.Font.SizeBi = 13

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rngBm = $d.Content
$rngBm.Find.Execute("Installation Guide", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $rngBm)

# ------------------------------------------------------------------
# 3. Stamp the font on the two now-empty paragraphs directly below the
#    title and below the authors/revision block, so their pPr carries
#    an explicit rPr (Helvetica Light).
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Font.NameAscii = "Helvetica Light"
$p3.Range.Font.NameOther = "Helvetica Light"

$p6 = $d.Paragraphs.Item(6)
$p6.Range.Font.NameAscii = "Helvetica Light"
$p6.Range.Font.NameOther = "Helvetica Light"

Write-Output "done"
